$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 720.1479260693817
$ws.Range("C2").Value = 11721.4951503772
$ws.Range("D2").Value = 9870.552512637782
